$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row holding "bocaioandoru12+2@gmail.com" (row 2); rows below shift up.
$ws.Rows.Item(2).Delete()

# The row delete does not retarget the worksheet's stale Hyperlinks collection
# (it still points old A2/A3/A4 at the original mailto addresses), so rebuild
# it to match the shifted data: A2 -> +3@gmail.com, A3 -> +4@gmail.com.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:bocaioandoru12+3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:bocaioandoru12+4@gmail.com")

# Re-adding a hyperlink stamps a fresh style; put the cells back on the
# original shared "Hyperlink" style so formatting matches pre-edit cells.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"

# Widen column A slightly to match the new layout (target stored width
# 27.44140625; the host quantizes ColumnWidth to 1/6-character pixel steps,
# so 26.667 is the input that lands on the nearest reachable stored width).
$ws.Columns.Item(1).ColumnWidth = 26.666666666666668

# Move the active selection to D4 (matches the saved selection state).
$ws.Range("D4").Select()
